$wb = $excel.ActiveWorkbook

# 1. Add a new "0" option to the barcode_offset lookup sheet (row 6).
#    Force text storage (so it lands in sharedStrings as "0" rather than
#    the number 0), then clear the formatting this introduces so the
#    cell doesn't pick up an extra quote-prefix style.
$barcodeOffsetSheet = $wb.Worksheets.Item("barcode_offset")
$cell = $barcodeOffsetSheet.Range("A6")
$cell.NumberFormat = "@"
$cell.Value = "0"
$cell.ClearFormats()

# 2. Extend the barcode_offset data validation on the ATACseq sheet so the
#    dropdown covers the new row (A1:A5 -> A1:A6).
$atacSeqSheet = $wb.Worksheets.Item("ATACseq")
$validationRange = $atacSeqSheet.Range("O2:O1001")
$validationRange.Validation.Formula1 = "='barcode_offset'!`$A`$1:`$A`$6"
$validationRange.Validation.IgnoreBlank = $true
$validationRange.Validation.ShowError = $true
$validationRange.Validation.AlertStyle = 1

# 3. Bump the pav:createdOn timestamp recorded on the .metadata sheet.
$metadataSheet = $wb.Worksheets.Item(".metadata")
$metadataSheet.Range("C2").Value = "2023-10-31T14:33:15-07:00"
